$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

for ($row = 8; $row -le 20; $row++) {
    $ws.Range("B$row").Value = 44743
    $ws.Range("C$row").Value = 44834
    $ws.Range("AD$row").Value = 44844
    $ws.Range("AE$row").Value = 44844
}

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 25
$ws.Range("AB8").Select()
